$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Text = "Legend: (match-mismatch-gap)`r"
